# "Generate Report for Handoff"
# The localization-status report was regenerated; for each language sheet
# (and the Overview sheet) the two rows describing
#   a65e57c5-1f67-41ac-9cbc-89d01f1280cc
#   b0c1d244-ec81-4b9e-975c-6d1bf13868a0
# swapped places (row 5 <-> row 6), and the handoff datetime for the file
# that is now on row 6 was bumped to a freshly generated timestamp.

$wb = $excel.ActiveWorkbook

function Set-DisplayForCell {
    param($ws, [string]$cellAddr, [string]$newText)

    $target = '$' + ($cellAddr -replace '(\d+)$', '$$$1')
    # cellAddr like "A5" -> target like "$A$5"
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            $hl.TextToDisplay = $newText
            return
        }
    }
}

# ---------- Overview sheet ----------
$wsOverview = $wb.Worksheets.Item("Overview")

$ov_a5 = $wsOverview.Range("A5").Value2
$ov_a6 = $wsOverview.Range("A6").Value2

$wsOverview.Range("A5").Value = $ov_a6
$wsOverview.Range("A6").Value = $ov_a5

Set-DisplayForCell $wsOverview "A5" $ov_a6
Set-DisplayForCell $wsOverview "A6" $ov_a5

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zh_a5 = $wsZh.Range("A5").Value2
$zh_a6 = $wsZh.Range("A6").Value2
$zh_c5 = $wsZh.Range("C5").Value2
$zh_c6 = $wsZh.Range("C6").Value2

$wsZh.Range("A5").Value = $zh_a6
$wsZh.Range("A6").Value = $zh_a5
$wsZh.Range("C5").Value = $zh_c6
$wsZh.Range("C6").Value = $zh_c5
$wsZh.Range("D5").Value = "2016-03-03 09:30:03"
$wsZh.Range("D6").Value = "2016-03-03 09:33:17"

Set-DisplayForCell $wsZh "A5" $zh_a6
Set-DisplayForCell $wsZh "A6" $zh_a5
Set-DisplayForCell $wsZh "C5" $zh_c6
Set-DisplayForCell $wsZh "C6" $zh_c5

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")

$de_a5 = $wsDe.Range("A5").Value2
$de_a6 = $wsDe.Range("A6").Value2
$de_c5 = $wsDe.Range("C5").Value2
$de_c6 = $wsDe.Range("C6").Value2

$wsDe.Range("A5").Value = $de_a6
$wsDe.Range("A6").Value = $de_a5
$wsDe.Range("C5").Value = $de_c6
$wsDe.Range("C6").Value = $de_c5
$wsDe.Range("D5").Value = "2016-03-03 09:30:19"
$wsDe.Range("D6").Value = "2016-03-03 09:33:29"

Set-DisplayForCell $wsDe "A5" $de_a6
Set-DisplayForCell $wsDe "A6" $de_a5
Set-DisplayForCell $wsDe "C5" $de_c6
Set-DisplayForCell $wsDe "C6" $de_c5

Write-Host "Done."
